$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002837
$ws.Range("H2").Value = 0.008510999999999999
$ws.Range("I2").Value = 0.00007108247730492929
$ws.Range("J2").Value = 0.00007108247730492929
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.01485317462584607
$ws.Range("P2").Value = 0.01485317462584607
$ws.Range("Q2").Value = 0.005162934308999999
$ws.Range("R2").Value = 0.046466408781
$ws.Range("S2").Value = 0.000001055800448247855
$ws.Range("T2").Value = 0.000001055800448247855
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002837
$ws.Range("H3").Value = 0.008510999999999999
$ws.Range("I3").Value = 0.00007108247730492929
$ws.Range("J3").Value = 0.00007108247730492929
$ws.Range("O3").Value = 0.726618572334523
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 0.25257118772
$ws.Range("R3").Value = 2.27314068948
$ws.Range("S3").Value = 0.00005164984817730885
$ws.Range("T3").Value = 0.00005164984817730886
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002837
$ws.Range("H4").Value = 0.008510999999999999
$ws.Range("I4").Value = 0.00007108247730492929
$ws.Range("J4").Value = 0.00007108247730492929
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.257333028084772
$ws.Range("P4").Value = 0.257333028084772
$ws.Range("Q4").Value = 0.08944845482566664
$ws.Range("R4").Value = 0.8050360934309998
$ws.Range("S4").Value = 0.00001829186912864454
$ws.Range("T4").Value = 0.00001829186912864454
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002837
$ws.Range("H5").Value = 0.008510999999999999
$ws.Range("I5").Value = 0.00007108247730492929
$ws.Range("J5").Value = 0.00007108247730492929
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.439328
$ws.Range("O5").Value = 0.001195224954858853
$ws.Range("P5").Value = 0.001195224954858853
$ws.Range("Q5").Value = 0.0004154578453333333
$ws.Range("R5").Value = 0.003739120608
$ws.Range("S5").Value = 0.00000008495955072803952
$ws.Range("T5").Value = 0.00000008495955072803953
$ws.Range("I6").Value = 0.3776915775490952
$ws.Range("J6").Value = 0.3776915775490952
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.01485317462584607
$ws.Range("P6").Value = 0.01485317462584607
$ws.Range("Q6").Value = 27.43287625702
$ws.Range("R6").Value = 246.89588631318
$ws.Range("S6").Value = 0.005609918956047993
$ws.Range("T6").Value = 0.005609918956047994
$ws.Range("I7").Value = 0.3776915775490952
$ws.Range("J7").Value = 0.3776915775490952
$ws.Range("O7").Value = 0.726618572334523
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("S7").Value = 0.2744377148614973
$ws.Range("T7").Value = 0.2744377148614974
$ws.Range("I8").Value = 0.3776915775490952
$ws.Range("J8").Value = 0.3776915775490952
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.257333028084772
$ws.Range("P8").Value = 0.257333028084772
$ws.Range("Q8").Value = 475.2778644377977
$ws.Range("R8").Value = 4277.500779940179
$ws.Range("S8").Value = 0.09719251733282318
$ws.Range("T8").Value = 0.09719251733282318
$ws.Range("I9").Value = 0.3776915775490952
$ws.Range("J9").Value = 0.3776915775490952
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.439328
$ws.Range("O9").Value = 0.001195224954858853
$ws.Range("P9").Value = 0.001195224954858853
$ws.Range("Q9").Value = 2.207505069582222
$ws.Range("R9").Value = 19.86754562624
$ws.Range("S9").Value = 0.0004514263987266861
$ws.Range("T9").Value = 0.0004514263987266862
$ws.Range("G10").Value = 1.581618666666667
$ws.Range("H10").Value = 4.744856
$ws.Range("I10").Value = 0.03962825977384063
$ws.Range("J10").Value = 0.03962825977384063
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.01485317462584607
$ws.Range("P10").Value = 0.01485317462584607
$ws.Range("Q10").Value = 2.878319801864
$ws.Range("R10").Value = 25.90487821677601
$ws.Range("S10").Value = 0.000588605462539246
$ws.Range("T10").Value = 0.0005886054625392461
$ws.Range("G11").Value = 1.581618666666667
$ws.Range("H11").Value = 4.744856
$ws.Range("I11").Value = 0.03962825977384063
$ws.Range("J11").Value = 0.03962825977384063
$ws.Range("O11").Value = 0.726618572334523
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 140.8076507437867
$ws.Range("R11").Value = 1267.26885669408
$ws.Range("S11").Value = 0.02879462954096968
$ws.Range("T11").Value = 0.02879462954096969
$ws.Range("G12").Value = 1.581618666666667
$ws.Range("H12").Value = 4.744856
$ws.Range("I12").Value = 0.03962825977384063
$ws.Range("J12").Value = 0.03962825977384063
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.257333028084772
$ws.Range("P12").Value = 0.257333028084772
$ws.Range("Q12").Value = 49.86723505701956
$ws.Range("R12").Value = 448.805115513176
$ws.Range("S12").Value = 0.01019766008533237
$ws.Range("T12").Value = 0.01019766008533237
$ws.Range("G13").Value = 1.581618666666667
$ws.Range("H13").Value = 4.744856
$ws.Range("I13").Value = 0.03962825977384063
$ws.Range("J13").Value = 0.03962825977384063
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.439328
$ws.Range("O13").Value = 0.001195224954858853
$ws.Range("P13").Value = 0.001195224954858853
$ws.Range("Q13").Value = 0.2316164551964444
$ws.Range("R13").Value = 2.084548096768
$ws.Range("S13").Value = 0.00004736468499932355
$ws.Range("T13").Value = 0.00004736468499932356
$ws.Range("G14").Value = 23.25273433333334
$ws.Range("H14").Value = 69.75820300000001
$ws.Range("I14").Value = 0.5826090801997593
$ws.Range("J14").Value = 0.5826090801997593
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.01485317462584607
$ws.Range("P14").Value = 0.01485317462584607
$ws.Range("Q14").Value = 42.31665134565701
$ws.Range("R14").Value = 380.8498621109131
$ws.Range("S14").Value = 0.008653594406810581
$ws.Range("T14").Value = 0.008653594406810582
$ws.Range("G15").Value = 23.25273433333334
$ws.Range("H15").Value = 69.75820300000001
$ws.Range("I15").Value = 0.5826090801997593
$ws.Range("J15").Value = 0.5826090801997593
$ws.Range("O15").Value = 0.726618572334523
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 2070.134201024894
$ws.Range("R15").Value = 18631.20780922404
$ws.Range("S15").Value = 0.4233345780838787
$ws.Range("T15").Value = 0.4233345780838788
$ws.Range("G16").Value = 23.25273433333334
$ws.Range("H16").Value = 69.75820300000001
$ws.Range("I16").Value = 0.5826090801997593
$ws.Range("J16").Value = 0.5826090801997593
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.257333028084772
$ws.Range("P16").Value = 0.257333028084772
$ws.Range("Q16").Value = 733.1410492028182
$ws.Range("R16").Value = 6598.269442825363
$ws.Range("S16").Value = 0.1499245587974879
$ws.Range("T16").Value = 0.1499245587974879
$ws.Range("G17").Value = 23.25273433333334
$ws.Range("H17").Value = 69.75820300000001
$ws.Range("I17").Value = 0.5826090801997593
$ws.Range("J17").Value = 0.5826090801997593
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.439328
$ws.Range("O17").Value = 0.001195224954858853
$ws.Range("P17").Value = 0.001195224954858853
$ws.Range("Q17").Value = 3.405192423064889
$ws.Range("R17").Value = 30.646731807584
$ws.Range("S17").Value = 0.0006963489115821149
$ws.Range("T17").Value = 0.0006963489115821151
